# 44344 consultation proposal line space
#
# 1. Split the "<Proposal Description>" paragraph (Keybody / italic style) into
#    two paragraphs: an empty leading paragraph, and a new paragraph that keeps
#    the pPr formatting plus the _GoBack bookmark (now collapsed/empty) and the
#    run with the placeholder text.
# 2. Add a <w:lastRenderedPageBreak/> marker in front of the "Yours sincerely"
#    run.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: "<Proposal Description>" paragraph split
# ---------------------------------------------------------------------------
$proposalPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*<Proposal Description>*") {
        $proposalPara = $p
    }
}

if ($proposalPara -ne $null) {
    $full = $proposalPara.Range

    $splitXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Keybody"/><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Keybody"/><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>&lt;Proposal Description&gt;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

    $full.InsertXML($splitXml)
}

# ---------------------------------------------------------------------------
# Part 2: "Yours sincerely" gets a lastRenderedPageBreak before the text
# ---------------------------------------------------------------------------
$signOffPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Yours sincerely*") {
        $signOffPara = $p
    }
}

if ($signOffPara -ne $null) {
    $full2 = $signOffPara.Range
    # exclude the trailing paragraph mark so the existing <w:p> is kept intact
    $content = $d.Range($full2.Start, $full2.End - 1)

    $pageBreakXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>Yours sincerely</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

    $content.InsertXML($pageBreakXml)
}
